$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1479.4166
$ws.Range("I43").Value = 700
$ws.Range("J43").Value = 1550.2727
$ws.Range("K43").Value = 700
$ws.Range("L43").Value = 1550.2727
$ws.Range("M43").Value = -631
$ws.Range("N43").Value = -1688.2727
# Row 62
$ws.Range("H62").Value = 2097.7693
$ws.Range("I62").Value = 1513
$ws.Range("J62").Value = 2463.25
$ws.Range("K62").Value = 1513
$ws.Range("L62").Value = 2463.25
$ws.Range("M62").Value = -889
$ws.Range("N62").Value = -3711.25
# Row 65
$ws.Range("H65").Value = 2097.7693
$ws.Range("I65").Value = 1513
$ws.Range("J65").Value = 2463.25
$ws.Range("K65").Value = 7565
$ws.Range("L65").Value = 12316.25
$ws.Range("M65").Value = -4445
$ws.Range("N65").Value = -18556.25
# Row 86
$ws.Range("H86").Value = 350034660
$ws.Range("I86").Value = 525050000
$ws.Range("J86").Value = 3999.5
$ws.Range("K86").Value = 525050000
$ws.Range("L86").Value = 3999.5
$ws.Range("M86").Value = -525048877
$ws.Range("N86").Value = -6245.5
# Row 89
$ws.Range("H89").Value = 350034660
$ws.Range("I89").Value = 525050000
$ws.Range("J89").Value = 3999.5
$ws.Range("K89").Value = 2625250000
$ws.Range("L89").Value = 19997.5
$ws.Range("M89").Value = -2625244384
$ws.Range("N89").Value = -31229.5
# Row 138
$ws.Range("H138").Value = 2322.2903
$ws.Range("I138").Value = 2366.4443
$ws.Range("J138").Value = 2311.6934
$ws.Range("K138").Value = 7099.3329
$ws.Range("L138").Value = 6935.0802
$ws.Range("M138").Value = -1959.3329
$ws.Range("N138").Value = -17215.0802

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1711.3793
$ws.Range("I74").Value = 1292.5758
$ws.Range("K74").Value = 1292.5758
$ws.Range("M74").Value = -418.5758000000001
# Row 77
$ws.Range("H77").Value = 1711.3793
$ws.Range("I77").Value = 1292.5758
$ws.Range("K77").Value = 6462.879000000001
$ws.Range("M77").Value = -2094.879000000001
# Row 122
$ws.Range("H122").Value = 102149.8
$ws.Range("I122").Value = 144471.14
$ws.Range("K122").Value = 433413.42
$ws.Range("M122").Value = -430963.42

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2606.9707
$ws.Range("I134").Value = 2382.2273
$ws.Range("K134").Value = 7146.6819
$ws.Range("M134").Value = -4611.6819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1140.8823
$ws.Range("I58").Value = 815.8095
$ws.Range("J58").Value = 1666
$ws.Range("K58").Value = 815.8095
$ws.Range("L58").Value = 1666
$ws.Range("M58").Value = -612.8095
$ws.Range("N58").Value = -2072
# Row 99
$ws.Range("H99").Value = 1883.1666
$ws.Range("I99").Value = 1833
$ws.Range("J99").Value = 1933.3334
$ws.Range("K99").Value = 1833
$ws.Range("L99").Value = 1933.3334
$ws.Range("M99").Value = -335
$ws.Range("N99").Value = -4929.3334
# Row 107
$ws.Range("H107").Value = 3290217.8
$ws.Range("I107").Value = 5208974
$ws.Range("J107").Value = 921.4286
$ws.Range("K107").Value = 5208974
$ws.Range("L107").Value = 921.4286
$ws.Range("M107").Value = -5207054
$ws.Range("N107").Value = -4761.4286
# Row 122
$ws.Range("H122").Value = 1513.1428
$ws.Range("I122").Value = 638
$ws.Range("K122").Value = 1914
$ws.Range("M122").Value = 536
# Row 126
$ws.Range("H126").Value = 1883.1666
$ws.Range("I126").Value = 1833
$ws.Range("J126").Value = 1933.3334
$ws.Range("K126").Value = 5499
$ws.Range("L126").Value = 5800.0002
$ws.Range("M126").Value = -3029
$ws.Range("N126").Value = -10740.0002
# Row 136
$ws.Range("H136").Value = 1140.8823
$ws.Range("I136").Value = 815.8095
$ws.Range("J136").Value = 1666
$ws.Range("K136").Value = 2447.4285
$ws.Range("L136").Value = 4998
$ws.Range("M136").Value = 102.5715
$ws.Range("N136").Value = -10098

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 246.21428
$ws.Range("J12").Value = 286.75
$ws.Range("L12").Value = 860.25
$ws.Range("N12").Value = -1206.25
# Row 17
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 250
$ws.Range("J17").Value = 750
$ws.Range("K17").Value = 750
$ws.Range("L17").Value = 2250
$ws.Range("M17").Value = -581
$ws.Range("N17").Value = -2588
# Row 39
$ws.Range("H39").Value = 1569.5775
$ws.Range("J39").Value = 1569.5775
$ws.Range("L39").Value = 4708.7325
$ws.Range("N39").Value = -5296.7325
# Row 64
$ws.Range("H64").Value = 2461.5386
$ws.Range("I64").Value = 800
$ws.Range("J64").Value = 2528
$ws.Range("K64").Value = 2400
$ws.Range("L64").Value = 7584
$ws.Range("M64").Value = -2130
$ws.Range("N64").Value = -8124
# Row 67
$ws.Range("H67").Value = 2461.5386
$ws.Range("I67").Value = 800
$ws.Range("J67").Value = 2528
$ws.Range("K67").Value = 2400
$ws.Range("L67").Value = 7584
$ws.Range("M67").Value = -1464
$ws.Range("N67").Value = -9456
# Row 82
$ws.Range("H82").Value = 2800
$ws.Range("I82").Value = 600
$ws.Range("J82").Value = 3166.6667
$ws.Range("K82").Value = 1800
$ws.Range("L82").Value = 9500.000100000001
$ws.Range("M82").Value = -1394
$ws.Range("N82").Value = -10312.0001
# Row 85
$ws.Range("H85").Value = 2800
$ws.Range("I85").Value = 600
$ws.Range("J85").Value = 3166.6667
$ws.Range("K85").Value = 1800
$ws.Range("L85").Value = 9500.000100000001
$ws.Range("M85").Value = -396
$ws.Range("N85").Value = -12308.0001
# Row 110
$ws.Range("H110").Value = 11359.823
$ws.Range("J110").Value = 12038.936
$ws.Range("L110").Value = 36116.808
$ws.Range("N110").Value = -44296.808
# Row 112
$ws.Range("H112").Value = 4907.6924
$ws.Range("I112").Value = 4400
$ws.Range("J112").Value = 5000
$ws.Range("K112").Value = 13200
$ws.Range("L112").Value = 15000
$ws.Range("M112").Value = -12092
$ws.Range("N112").Value = -17216

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 244739.08
$ws.Range("I11").Value = 222955.56
$ws.Range("J11").Value = 293752
$ws.Range("K11").Value = 222955.56
$ws.Range("L11").Value = 293752
$ws.Range("M11").Value = -222816.56
$ws.Range("N11").Value = -294030
# Row 19
$ws.Range("H19").Value = 46668.668
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 46668.668
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 46668.668
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -47244.668
# Row 113
$ws.Range("H113").Value = 251112.5
$ws.Range("I113").Value = 500375
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 500375
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = -498205
$ws.Range("N113").Value = -6190
# Row 126
$ws.Range("H126").Value = 2084.1667
$ws.Range("I126").Value = 1916.75
$ws.Range("J126").Value = 2419
$ws.Range("K126").Value = 5750.25
$ws.Range("L126").Value = 7257
$ws.Range("M126").Value = -3280.25
$ws.Range("N126").Value = -12197

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1464.091
$ws.Range("I7").Value = 1110
$ws.Range("K7").Value = 1110
$ws.Range("M7").Value = -998
# Row 32
$ws.Range("H32").Value = 789.5
$ws.Range("I32").Value = 789.5
$ws.Range("K32").Value = 789.5
$ws.Range("M32").Value = -472.5
# Row 126
$ws.Range("H126").Value = 1464.091
$ws.Range("I126").Value = 1110
$ws.Range("K126").Value = 3330
$ws.Range("M126").Value = -860

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1012
$ws.Range("I113").Value = 1012
$ws.Range("K113").Value = 3036
$ws.Range("M113").Value = -866
# Row 122
$ws.Range("H122").Value = 5004
$ws.Range("I122").Value = 5004
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15012
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12562
$ws.Range("N122").ClearContents()
